# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" palette (used by the notes master)
#   ppt/theme/theme2.xml -> "Integral" palette      (used by the slide master / design)
#
# The target revision swaps the two palettes: the design's theme (serialized to
# theme2.xml) becomes the "Office Theme" colors, while the notes-master theme
# (theme1.xml) becomes the "Integral" colors. This host exposes a single
# mutable Theme object (reached via SlideMaster.Theme / NotesMaster.Theme /
# Slide.ThemeColorScheme, etc. - they all resolve to the same design), so we
# drive the reachable half of the swap: recolor the presentation's theme with
# the "Office Theme" values.

function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$cs = $theme.ThemeColorScheme

# Target palette ("Office Theme"), in the fixed 12-slot theme-color order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink.
$cs.Colors(1).RGB  = RGBVal 0x00 0x00 0x00   # dk1
$cs.Colors(2).RGB  = RGBVal 0xFF 0xFF 0xFF   # lt1
$cs.Colors(3).RGB  = RGBVal 0x44 0x54 0x6A   # dk2
$cs.Colors(4).RGB  = RGBVal 0xE7 0xE6 0xE6   # lt2
$cs.Colors(5).RGB  = RGBVal 0x5B 0x9B 0xD5   # accent1
$cs.Colors(6).RGB  = RGBVal 0xED 0x7D 0x31   # accent2
$cs.Colors(7).RGB  = RGBVal 0xA5 0xA5 0xA5   # accent3
$cs.Colors(8).RGB  = RGBVal 0xFF 0xC0 0x00   # accent4
$cs.Colors(9).RGB  = RGBVal 0x44 0x72 0xC4   # accent5
$cs.Colors(10).RGB = RGBVal 0x70 0xAD 0x47   # accent6
$cs.Colors(11).RGB = RGBVal 0x05 0x63 0xC1   # hlink
$cs.Colors(12).RGB = RGBVal 0x95 0x4F 0x72   # folHlink
